# Lesson1.xlsx update: add "Comment khi push" guidance block (rows 22-25),
# move the lesson-number cell from D2 to E2 (and update the two formulas
# that reference it), and record the cursor position on Sheet2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lesson1")

# ---------------------------------------------------------------------
# 1) Move the lesson number from D2 to E2 (copy format, move value, then
#    clear the old cell entirely so no stray <c r="D2"/> remains).
# ---------------------------------------------------------------------
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value2 = 1
$ws.Range("D2").Clear()

# ---------------------------------------------------------------------
# 2) Row 20 formula now points at E2 instead of D2.
# ---------------------------------------------------------------------
$ws.Range("D20").Formula = '="Lesson-"&E2&"-{tên user}"'

# ---------------------------------------------------------------------
# 3) Row 21: the example formula shifts from C21 to D21 and also now
#    points at E2. Copy C21's formatting onto D21 first, write the new
#    formula, then clear C21 completely.
# ---------------------------------------------------------------------
$ws.Range("C21").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("D21").Formula = '="vd: Lesson-"&E2&"-an"'
$ws.Range("C21").Clear()

# ---------------------------------------------------------------------
# 4) New "Comment khi push" block, rows 22-25, styled like the existing
#    red "Note" block above it (copy format from C18, which already
#    carries that style, onto the new C22:L25 range).
# ---------------------------------------------------------------------
$ws.Range("C18").Copy()
$ws.Range("C22:L25").PasteSpecial(-4122)

# Write the new text in the order that matches the target shared-string
# table layout (strings are registered in first-write order).
$ws.Range("C22").Value2 = "Comment khi push"
$ws.Range("D23").Value2 = "#{số issue}: comment"
$ws.Range("E25").Value2 = "#1: hoàn thành lesson 1"
$ws.Range("D24").Value2 = "vd: Khi hoàn thành lesson 1 ở https://github.com/VHEC-NAMTT/asp_net_training/issues/1 khi push sẽ ghi"

# ---------------------------------------------------------------------
# 5) Record the last-used cursor position on Sheet2 (B22), then return
#    to Lesson1 so it stays the active sheet/tab.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("B22").Select()
$ws.Activate()
